# Add 7 new merge-result columns (Action .. EachCost_y) to the right of the
# existing Date/Stock_Id/ProfitPercent/ProfitMoney table, matching the
# header style already used by the existing header row (bold, bordered,
# centered, top-aligned). The new data rows (2-9) are left blank, since the
# source data update only populated these extra columns for the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @("Action", "PositionSize_x", "Price_x", "EachCost_x", "PositionSize_y", "Price_y", "EachCost_y")

# Copy the formatting of the last existing header cell (D1) onto each new
# header cell (E1:K1) so they reuse the same cell style as A1:D1, then set
# the text.
$ws.Range("D1").Copy() | Out-Null
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $ws.Cells.Item(1, 5 + $i)
    $cell.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $cell.Value = $headers[$i]
}

Write-Output "done"
